$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'roxy leggings',
    'rugby 16',
    'rugby 2',
    'rugby apparel',
    'rugby apparel for women',
    'rugby bag',
    'rugby bags',
    'rugby ball all blacks',
    'rugby ball womens',
    'rugby challenge 3',
    'rugby dress women',
    'rugby fitness',
    'rugby gear',
    'rugby gear for men',
    'rugby jackets for women',
    'rugby knee brace',
    'rugby live',
    'rugby pants',
    'rugby patches',
    'rugby pro',
    'rugby short',
    'rugby shorts',
    'rugby shorts women',
    'rugby skins',
    'rugby tank',
    'rugby team shirts',
    'rugby training equipment',
    'rugby women',
    'rugby womens shorts',
    'run apparel',
    'run compression',
    'run cow',
    'run leggings for women',
    'run pants women',
    'run sausage run game',
    'run shorts for women',
    'run shorts women',
    'run tights women',
    'run track',
    'run training equipment',
    'runing equipment',
    'runner apparel',
    'runner bike',
    'runner chafing',
    'runner compression pants',
    'runner equipment',
    'runner gear',
    'runner gift basket',
    'runner jacket womens',
    'runner knee brace',
    'runner knee support',
    'runner leggings women',
    'runner protection for women ring',
    'runner recovery gear',
    'runner reflective gear',
    'runner ring',
    'runner rings for women',
    'runner shorts women',
    'runner supplement',
    'runner tanks for women',
    'runner waist band',
    'runner waist pack',
    'runners basket',
    'runners clothes',
    'runners clothes cold weather',
    'runners gear for men',
    'runners gear for women',
    'runners high',
    'runners jacket women',
    'runners knee',
    'runners knee brace woman',
    'runners knee brace women',
    'runners knee support',
    'runners knee tape',
    'runners leggings',
    'runners leggings women',
    'runners legs',
    'runners neck warmer women',
    'runners pack women',
    'runners reflective',
    'runners reflective gear',
    'runners shorts women',
    'runners shorts womens',
    'runners skirt',
    'runners tank tops for women',
    'runners tape knee',
    'runners underwear women',
    'runners vest with pockets',
    'runners winter jacket women',
    'running 2xu',
    'running 80 20',
    'running apparel',
    'running apparel for women',
    'running apparel women',
    'running apparell',
    'running apperal',
    'running back football gear',
    'running band knee',
    'running butter',
    'running capri'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
